$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") contains a date serial that should be bumped
# from 45177 (2023-09-08) to 45178 (2023-09-09) for every data row.
$lastRow = 422
$ws.Range("C2:C$lastRow").Value = 45178
